$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title (October 2016 -> November 2016)
$ws.Range("A1").Value = "Table 4.16. Receipts and Quality of Coal by Rank Delivered for Electricity Generation: Independent Power Producers by State, November 2016"

# Row 5 - New England
$ws.Range("B5").Value = 90
$ws.Range("C5").Value = 0.64
$ws.Range("D5").Value = 9.4

# Row 7 - Maine
$ws.Range("B7").Value = 5

# Row 8 - Massachusetts
$ws.Range("B8").Value = 85
$ws.Range("C8").Value = 0.63
$ws.Range("D8").Value = 9.5

# Row 12 - Middle Atlantic
$ws.Range("B12").Value = 1394
$ws.Range("C12").Value = 3.3
$ws.Range("D12").Value = 8.9

# Row 13 - New Jersey
$ws.Range("B13").Value = 52
$ws.Range("C13").Value = 1.77
$ws.Range("D13").Value = 7.2

# Row 14 - New York
$ws.Range("B14").Value = 44
$ws.Range("C14").Value = 2.87
$ws.Range("D14").Value = 8.2

# Row 15 - Pennsylvania
$ws.Range("B15").Value = 1297
$ws.Range("C15").Value = 3.37
$ws.Range("D15").Value = 9

# Row 16 - East North Central
$ws.Range("B16").Value = 1812
$ws.Range("C16").Value = 3.39
$ws.Range("D16").Value = 11.2
$ws.Range("E16").Value = 2323
$ws.Range("F16").Value = 0.21
$ws.Range("G16").Value = 4.6

# Row 17 - Illinois
$ws.Range("B17").Value = 267
$ws.Range("C17").Value = 3.97
$ws.Range("D17").Value = 27.7
$ws.Range("E17").Value = 2323
$ws.Range("F17").Value = 0.21
$ws.Range("G17").Value = 4.6

# Row 18 - Indiana
$ws.Range("B18").Value = 101
$ws.Range("C18").Value = 3.7
$ws.Range("D18").Value = 9.2

# Row 19 - Michigan
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = 0.48
$ws.Range("D19").Value = 5.8

# Row 20 - Ohio
$ws.Range("B20").Value = 1426
$ws.Range("C20").Value = 3.32
$ws.Range("D20").Value = 9.2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = "--"
$ws.Range("G20").Value = "--"

# Row 30 - West North Central
$ws.Range("B30").Value = 1048
$ws.Range("D30").Value = 10.5
$ws.Range("E30").Value = 27
$ws.Range("F30").Value = 0.2
$ws.Range("G30").Value = 4.6

# Row 31 - Iowa
$ws.Range("B31").Value = 14
$ws.Range("C31").Value = 2.41
$ws.Range("D31").Value = 7.5

# Row 33 - Minnesota
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = "--"
$ws.Range("D33").Value = "--"

# Row 35 - Missouri
$ws.Range("B35").Value = 588
$ws.Range("C35").Value = 2.59
$ws.Range("D35").Value = 9
$ws.Range("E35").Value = 27
$ws.Range("F35").Value = 0.2
$ws.Range("G35").Value = 4.6

# Row 38 - North Dakota
$ws.Range("B38").Value = 67
$ws.Range("C38").Value = 0.81
$ws.Range("D38").Value = 9

# Row 39 - South Dakota
$ws.Range("B39").Value = 379
$ws.Range("C39").Value = 3.87
$ws.Range("D39").Value = 13.3

# Row 40 - South Atlantic
$ws.Range("H40").Value = 261
$ws.Range("I40").Value = 0.48
$ws.Range("J40").Value = 14.1

# Row 43 - Florida
$ws.Range("H43").Value = 261
$ws.Range("I43").Value = 0.48
$ws.Range("J43").Value = 14.1

# Row 45 - East South Central
$ws.Range("B45").Value = 33
$ws.Range("C45").Value = 1.89
$ws.Range("D45").Value = 25.8
$ws.Range("E45").Value = 3566
$ws.Range("F45").Value = 0.3
$ws.Range("G45").Value = 5.3
$ws.Range("H45").Value = 2448
$ws.Range("I45").Value = 0.93
$ws.Range("J45").Value = 15.7

# Row 46 - Alabama
$ws.Range("F46").Value = 0.23

# Row 47 - Kentucky
$ws.Range("E47").Value = 230
$ws.Range("G47").Value = 5

# Row 48 - Mississippi
$ws.Range("B48").Value = 33
$ws.Range("C48").Value = 1.89
$ws.Range("D48").Value = 25.8
$ws.Range("E48").Value = 56
$ws.Range("F48").Value = 0.4
$ws.Range("G48").Value = 5

# Row 49 - Tennessee
$ws.Range("E49").Value = 3064
$ws.Range("F49").Value = 0.3
$ws.Range("G49").Value = 5.3
$ws.Range("H49").Value = 2448
$ws.Range("I49").Value = 0.93
$ws.Range("J49").Value = 15.7

# Row 50 - West South Central
$ws.Range("E50").Value = 843
$ws.Range("F50").Value = 0.68
$ws.Range("G50").Value = 9.9

# Row 54 - Oklahoma
$ws.Range("E54").Value = 803
$ws.Range("F54").Value = 0.7
$ws.Range("G54").Value = 10.2

# Row 55 - Texas
$ws.Range("E55").Value = 40
$ws.Range("F55").Value = 0.32
$ws.Range("G55").Value = 5.4

# Row 59 - Mountain
$ws.Range("E59").Value = 475
$ws.Range("F59").Value = 0.43
$ws.Range("G59").Value = 8.8

# Row 62 - Utah
$ws.Range("E62").Value = 475
$ws.Range("F62").Value = 0.43
$ws.Range("G62").Value = 8.8

# Row 63 - Pacific Contiguous
$ws.Range("E63").Value = 120
$ws.Range("F63").Value = 0.19
$ws.Range("G63").Value = 3.8

# Row 65 - Washington
$ws.Range("E65").Value = 120
$ws.Range("F65").Value = 0.19
$ws.Range("G65").Value = 3.8

# Row 66 - U.S. Total
$ws.Range("B66").Value = 4377
$ws.Range("D66").Value = 10.3
$ws.Range("E66").Value = 7354
$ws.Range("F66").Value = 0.32
$ws.Range("H66").Value = 2709
$ws.Range("I66").Value = 0.9
$ws.Range("J66").Value = 15.6
